# Updates cryptos list data (prices & 1h volume %) per upstream scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.661.50"
$ws.Range("E2").Value = "  +8.74%  "

# Row 3
$ws.Range("D3").Value = "3.645.36"
$ws.Range("E3").Value = "  +8.98%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'596.96"
$ws.Range("E5").Value = "  +7.34%  "

# Row 6
$ws.Range("D6").Value = "'193.65"

# Row 7
$ws.Range("E7").Value = "  +4.19%  "

# Row 8
$ws.Range("D8").Value = "3.605.48"
$ws.Range("E8").Value = "  +8.08%  "

# Row 9
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  +7.79%  "

# Row 11
$ws.Range("E11").Value = "  +6.58%  "

# Row 12
$ws.Range("D12").Value = "'58.14"
$ws.Range("E12").Value = "  +10.69%  "

# Row 13
$ws.Range("D13").Value = "'0.0000305"
$ws.Range("E13").Value = "  +11.20%  "

# Row 14
$ws.Range("D14").Value = "'9.81"
$ws.Range("E14").Value = "  +8.09%  "

# Row 15
$ws.Range("D15").Value = "4.218.45"
$ws.Range("E15").Value = "  +8.65%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.639.28"
$ws.Range("E16").Value = "  +8.56%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.47"
$ws.Range("E17").Value = "  +7.98%  "

# Row 18
$ws.Range("D18").Value = "70.562.31"
$ws.Range("E18").Value = "  +8.62%  "

# Row 19
$ws.Range("D19").Value = "'12.70"
$ws.Range("E19").Value = "  +8.99%  "

# Row 20
$ws.Range("E20").Value = "  +1.54%  "

# Row 21
$ws.Range("E21").Value = "  +7.60%  "

# Row 22
$ws.Range("D22").Value = "'499.64"
$ws.Range("E22").Value = "  +5.12%  "

# Row 23
$ws.Range("D23").Value = "'5.54"
$ws.Range("E23").Value = "  +13.42%  "

# Row 24
$ws.Range("D24").Value = "'16.81"
$ws.Range("E24").Value = "  +20.05%  "

# Row 25
$ws.Range("E25").Value = "  +10.73%  "

# Row 26
$ws.Range("D26").Value = "'91.39"
$ws.Range("E26").Value = "  +2.10%  "

# Row 27
$ws.Range("D27").Value = "'3.13"
$ws.Range("E27").Value = "  +9.47%  "

# Row 28
$ws.Range("D28").Value = "'11.33"
$ws.Range("E28").Value = "  +8.42%  "

# Row 29
$ws.Range("D29").Value = "'9.44"
$ws.Range("E29").Value = "  +10.12%  "

# Row 30
$ws.Range("D30").Value = "'32.65"
$ws.Range("E30").Value = "  +5.77%  "

# Row 31
$ws.Range("D31").Value = "'7.56"
$ws.Range("E31").Value = "  +17.45%  "

# Row 32
$ws.Range("D32").Value = "'12.24"
$ws.Range("E32").Value = "  +8.36%  "

# Row 33
$ws.Range("D33").Value = "'619.21"
$ws.Range("E33").Value = "  +8.42%  "

# Row 34
$ws.Range("D34").Value = "'65.86"
$ws.Range("E34").Value = "  +6.83%  "

# Row 35
$ws.Range("E35").Value = "  +9.98%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0842"
$ws.Range("E36").Value = "  +15.82%  "

# Row 37
$ws.Range("E37").Value = "  +6.56%  "

# Row 38
$ws.Range("D38").Value = "'3.73"
$ws.Range("E38").Value = "  +3.26%  "

# Row 39
$ws.Range("D39").Value = "'38.25"
$ws.Range("E39").Value = "  +8.91%  "

# Row 40
$ws.Range("D40").Value = "'0.404"
$ws.Range("E40").Value = "  +9.84%  "

# Row 41
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("D42").Value = "3.343.12"
$ws.Range("E42").Value = "  +8.70%  "

# Row 43
$ws.Range("E43").Value = "  +12.29%  "

# Row 44
$ws.Range("E44").Value = "  +9.30%  "

# Row 45
$ws.Range("D45").Value = "'2.69"
$ws.Range("E45").Value = "  +12.52%  "

# Row 46
$ws.Range("D46").Value = "'2.80"
$ws.Range("E46").Value = "  +18.31%  "

# Row 47
$ws.Range("E47").Value = "  +3.68%  "

# Row 48
$ws.Range("E48").Value = "  +4.12%  "

# Row 49
$ws.Range("D49").Value = "'9.11"
$ws.Range("E49").Value = "  +9.35%  "

# Row 50
$ws.Range("E50").Value = "  +7.24%  "

# Row 51
$ws.Range("E51").Value = "  -0.13%  "
